$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "iaest-measure:nivel-estudios-detalle"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "iaest-measure:nivel-estudios-agregado"

$ws.Range("B3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"

$ws.Range("B4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"

$ws.Range("A5:J5").EntireRow.Delete()
